$d = $word.ActiveDocument

# Fix mojibake Euro sign: "â‚¬" (the UTF-8 bytes for "€" were mis-decoded as
# Latin-1/cp1252 and re-encoded) should be the correct Euro sign "€" in the
# table header text:
#   "(million â‚¬ 2000 - 2006)"  ->  "(million € 2000 - 2006)"

$old = "million " + [char]0x00E2 + [char]0x201A + [char]0x00AC
$new = "million " + [char]0x20AC

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $new, 2)
